# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 23:14"

# Helper table of updated values: Row => @{ Col => NewValue }
$updates = @{
    4   = @{ B=4755609; C=49720;  D=2354194; E=2243693;          G=975; H=157722 }
    5   = @{ B=2707877; C=41579;            E=730263;            G=995; H=93563  }
    6   = @{ B=1751919; C=54865;  D=1146879; E=567637 }
    8   = @{ B=503290;  C=10107;  D=342461;  E=152676;            G=148; H=8153  }
    25  = @{ B=116582;  C=270;    D=101390;  E=6251;              G=6;   H=8941  }
    29  = @{ B=94316;   C=238;    D=41137;   E=48345;             G=29;  H=4834  }
    36  = @{ B=72218;   C=1248;   D=45102;   E=26590;             G=14;  H=526   }
    57  = @{           D=31300;  E=2131 }
    107 = @{ B=4186;    C=108;    D=1914;    E=2152;              G=6;   H=120   }
    138 = @{ B=1552;    C=17;     D=1217;    E=284;               G=1;   H=51    }
    149 = @{ B=958;     C=17;     D=647;     E=292 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
